# add: reset password and search
# (applies the jurusan_template.xlsx content/formatting refresh: renamed
#  "Nama Jurusan" values get a degree-level prefix, new D3 row appended,
#  whole table restyled to Times New Roman 12 with centered alignment.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Update existing "Nama Jurusan" text values (add S1/S2 prefixes) ----
$ws.Range("B2").Value = "S1 Akuakultur"
$ws.Range("B3").Value = "S1 Akuntansi"
$ws.Range("B4").Value = "S1 Hukum"
$ws.Range("B5").Value = "S1 Ilmu Hubungan Internasional"
$ws.Range("B6").Value = "S1 Ilmu Komunikasi"
$ws.Range("B7").Value = "S2 Magister Manajemen"
$ws.Range("B8").Value = "S1 Manajemen"
$ws.Range("B9").Value = "S1 Pemanfaatan Sumber Daya Perikanan"
$ws.Range("B10").Value = "S1 Sistem Informasi"
$ws.Range("B11").Value = "S1 Teknik Informatika"
$ws.Range("B12").Value = "S1 Teknik Lingkungan"

# ---- 2. Append the new D3 Manajemen Informatika row ----
$ws.Range("A13").Value = "57401"
$ws.Range("B13").Value = "D3 Manajemen Informatika"

# ---- 3. Re-font the whole table (Times New Roman 12) ----
$ws.Range("A1:B13").Font.Name = "Times New Roman"
$ws.Range("A1:B13").Font.Size = 12
$ws.Range("K9").Font.Name = "Times New Roman"
$ws.Range("K9").Font.Size = 12

# header row keeps its bold weight
$ws.Range("A1:B1").Font.Bold = $true

# ---- 4. Center alignment (horizontal + vertical) across the table ----
$ws.Range("A1:B13").HorizontalAlignment = -4108
$ws.Range("A1:B13").VerticalAlignment = -4108

# ---- 5. Row heights / column widths ----
$ws.Rows("1:13").RowHeight = 15.75
$ws.Columns.Item(1).ColumnWidth = 14.42578125
$ws.Columns.Item(2).ColumnWidth = 34.42578125

# ---- 6. Selection moves to B5 ----
$ws.Range("B5").Select()
